# DLAD Part 30 - TOC link fix:
# The Table of Contents entry for "30.201-5" hyperlinks to the bookmark
# anchor "P30_201_5", but that bookmark never existed in the body heading.
# Add the missing bookmark around the "30.201-5 " portion of the
# "30.201-5 Waiver." heading so the TOC hyperlink resolves correctly.

$d = $word.ActiveDocument

# Locate the body heading paragraph "30.201-5 Waiver." (the actual
# section heading, not the earlier Table-of-Contents line that merely
# links to it).
$target = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "30.201-5 Waiver.") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the '30.201-5 Waiver.' heading paragraph"
}

$full = $target.Range
$headStart = $full.Start
$bmRange = $d.Range($headStart, $headStart + 9)

# Sanity check: the range must cover exactly "30.201-5 " (with trailing
# space) before we bookmark it.
if ($bmRange.Text -ne "30.201-5 ") {
    throw "Unexpected range text: [$($bmRange.Text)]"
}

$d.Bookmarks.Add("P30_201_5", $bmRange)
